$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update Total row (B12): 36 -> 60
$ws.Range("B12").Value = 60

# Update the score display (E12): "35/84" -> "60/140"
$ws.Range("E12").Value = "60/140"
